$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 5000
$ws.Range("J45").Value = 5000
$ws.Range("L45").Value = 15000
$ws.Range("N45").Value = -15384
$ws.Range("H51").Value = 31396.6
$ws.Range("J51").Value = 19246
$ws.Range("L51").Value = 19246
$ws.Range("N51").Value = -20214
$ws.Range("H53").Value = 3139.9524
$ws.Range("I53").Value = 4513
$ws.Range("J53").Value = 2453.4285
$ws.Range("K53").Value = 4513
$ws.Range("L53").Value = 2453.4285
$ws.Range("M53").Value = -3876
$ws.Range("N53").Value = -3727.4285
$ws.Range("H116").Value = 62512500
$ws.Range("I116").Value = 250000000
$ws.Range("K116").Value = 250000000
$ws.Range("M116").Value = -249996558
$ws.Range("H121").Value = 19574.75
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 19574.75
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 58724.25
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -62218.25
$ws.Range("H132").Value = 1949.125
$ws.Range("I132").Value = 1949.125
$ws.Range("K132").Value = 5847.375
$ws.Range("M132").Value = -3317.375
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 7682.0586
$ws.Range("I137").Value = 6333.5
$ws.Range("J137").Value = 8417.637000000001
$ws.Range("K137").Value = 19000.5
$ws.Range("L137").Value = 25252.911
$ws.Range("M137").Value = -16450.5
$ws.Range("N137").Value = -30352.911
$ws.Range("H138").Value = 9815415
$ws.Range("I138").Value = 3510.5
$ws.Range("J138").Value = 14432782
$ws.Range("K138").Value = 10531.5
$ws.Range("L138").Value = 43298346
$ws.Range("M138").Value = -5391.5
$ws.Range("N138").Value = -43308626

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2504528.5
$ws.Range("I32").Value = 2504528.5
$ws.Range("K32").Value = 2504528.5
$ws.Range("M32").Value = -2504241.5
$ws.Range("H61").Value = 52639470
$ws.Range("I61").Value = 1451.875
$ws.Range("K61").Value = 1451.875
$ws.Range("M61").Value = -1239.875
$ws.Range("H74").Value = 137834.42
$ws.Range("I74").Value = 534333
$ws.Range("K74").Value = 534333
$ws.Range("M74").Value = -533459
$ws.Range("H77").Value = 137834.42
$ws.Range("I77").Value = 534333
$ws.Range("K77").Value = 2671665
$ws.Range("M77").Value = -2667297
$ws.Range("H107").Value = 59338
$ws.Range("J107").Value = 59338
$ws.Range("L107").Value = 59338
$ws.Range("N107").Value = -67018
$ws.Range("H110").Value = 27778476
$ws.Range("I110").Value = 745.7
$ws.Range("J110").Value = 166667120
$ws.Range("K110").Value = 745.7
$ws.Range("L110").Value = 166667120
$ws.Range("M110").Value = 1299.3
$ws.Range("N110").Value = -166671210
$ws.Range("H122").Value = 13094.2
$ws.Range("I122").Value = 19255.615
$ws.Range("K122").Value = 57766.845
$ws.Range("M122").Value = -55316.845
$ws.Range("H136").Value = 52639470
$ws.Range("I136").Value = 1451.875
$ws.Range("K136").Value = 4355.625
$ws.Range("M136").Value = -1805.625

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 283.13333
$ws.Range("I22").Value = 277
$ws.Range("K22").Value = 277
$ws.Range("M22").Value = -104
$ws.Range("H94").Value = 6513.8
$ws.Range("I94").Value = 6282
$ws.Range("K94").Value = 6282
$ws.Range("M94").Value = -5831
$ws.Range("H105").Value = 3218.44
$ws.Range("I105").Value = 2279.0833
$ws.Range("K105").Value = 2279.0833
$ws.Range("M105").Value = -532.0832999999998
$ws.Range("H107").Value = 29610792
$ws.Range("I107").Value = 36295988
$ws.Range("J107").Value = 4923.143
$ws.Range("K107").Value = 36295988
$ws.Range("L107").Value = 4923.143
$ws.Range("M107").Value = -36294068
$ws.Range("N107").Value = -8763.143

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5769.339
$ws.Range("I31").Value = 2318.08
$ws.Range("J31").Value = 8307.029
$ws.Range("K31").Value = 2318.08
$ws.Range("L31").Value = 8307.029
$ws.Range("M31").Value = -2023.08
$ws.Range("N31").Value = -8897.029
$ws.Range("H34").Value = 5769.339
$ws.Range("I34").Value = 2318.08
$ws.Range("J34").Value = 8307.029
$ws.Range("K34").Value = 2318.08
$ws.Range("L34").Value = 8307.029
$ws.Range("M34").Value = -2116.08
$ws.Range("N34").Value = -8711.029
$ws.Range("H105").Value = 6495636
$ws.Range("I105").Value = 7937834
$ws.Range("K105").Value = 7937834
$ws.Range("M105").Value = -7936087
$ws.Range("H107").Value = 2607.0833
$ws.Range("I107").Value = 2496.7273
$ws.Range("K107").Value = 2496.7273
$ws.Range("M107").Value = -576.7273
$ws.Range("H122").Value = 38479.68
$ws.Range("I122").Value = 2108.111
$ws.Range("J122").Value = 103948.5
$ws.Range("K122").Value = 6324.333
$ws.Range("L122").Value = 311845.5
$ws.Range("M122").Value = -3874.333
$ws.Range("N122").Value = -316745.5
$ws.Range("H140").Value = 49993
$ws.Range("J140").Value = 49993
$ws.Range("L140").Value = 49993
$ws.Range("N140").Value = -60353

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5007239
$ws.Range("J55").Value = 5890518
$ws.Range("L55").Value = 17671554
$ws.Range("N55").Value = -17671908
$ws.Range("H122").Value = 5657486.5
$ws.Range("I122").Value = 9428143
$ws.Range("K122").Value = 84853287
$ws.Range("M122").Value = -84850837
$ws.Range("H137").Value = 145169.42
$ws.Range("J137").Value = 203788.2
$ws.Range("L137").Value = 611364.6000000001
$ws.Range("N137").Value = -621564.6000000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 579370.6
$ws.Range("I70").Value = 806319
$ws.Range("J70").Value = 11999.75
$ws.Range("K70").Value = 806319
$ws.Range("L70").Value = 11999.75
$ws.Range("M70").Value = -806049
$ws.Range("N70").Value = -12539.75
$ws.Range("H73").Value = 579370.6
$ws.Range("I73").Value = 806319
$ws.Range("J73").Value = 11999.75
$ws.Range("K73").Value = 806319
$ws.Range("L73").Value = 11999.75
$ws.Range("M73").Value = -805383
$ws.Range("N73").Value = -13871.75
$ws.Range("H122").Value = 6040047
$ws.Range("I122").Value = 9057228
$ws.Range("K122").Value = 27171684
$ws.Range("M122").Value = -27169234

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4373.35
$ws.Range("I7").Value = 2698
$ws.Range("K7").Value = 2698
$ws.Range("M7").Value = -2586
$ws.Range("H22").Value = 14034.444
$ws.Range("I22").Value = 862.25
$ws.Range("J22").Value = 24572.2
$ws.Range("K22").Value = 862.25
$ws.Range("L22").Value = 24572.2
$ws.Range("M22").Value = -567.25
$ws.Range("N22").Value = -25162.2
$ws.Range("H27").Value = 14034.444
$ws.Range("I27").Value = 862.25
$ws.Range("J27").Value = 24572.2
$ws.Range("K27").Value = 862.25
$ws.Range("L27").Value = 24572.2
$ws.Range("M27").Value = -755.25
$ws.Range("N27").Value = -24786.2
$ws.Range("H40").Value = 4178.0605
$ws.Range("J40").Value = 5991.9165
$ws.Range("L40").Value = 5991.9165
$ws.Range("N40").Value = -6263.9165
$ws.Range("H55").Value = 392.13635
$ws.Range("I55").Value = 29.285715
$ws.Range("K55").Value = 29.285715
$ws.Range("M55").Value = 143.714285
$ws.Range("H94").Value = 35000
$ws.Range("J94").Value = 35000
$ws.Range("L94").Value = 35000
$ws.Range("N94").Value = -36352
$ws.Range("H122").Value = 3517
$ws.Range("I122").Value = 2643.45
$ws.Range("K122").Value = 7930.349999999999
$ws.Range("M122").Value = -5480.349999999999
$ws.Range("H126").Value = 4373.35
$ws.Range("I126").Value = 2698
$ws.Range("K126").Value = 8094
$ws.Range("M126").Value = -5624
$ws.Range("H136").Value = 9345.147000000001
$ws.Range("I136").Value = 3327.3845
$ws.Range("J136").Value = 13070.429
$ws.Range("K136").Value = 9982.1535
$ws.Range("L136").Value = 39211.287
$ws.Range("M136").Value = -7432.1535
$ws.Range("N136").Value = -44311.287

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 21105952
$ws.Range("I81").Value = 1558026
$ws.Range("K81").Value = 3116052
$ws.Range("M81").Value = -3114991
$ws.Range("H84").Value = 21105952
$ws.Range("I84").Value = 1558026
$ws.Range("K84").Value = 15580260
$ws.Range("M84").Value = -15574956
$ws.Range("H122").Value = 14005425
$ws.Range("I122").Value = 17384818
$ws.Range("K122").Value = 52154454
$ws.Range("M122").Value = -52152004
